# Applies the scheduled-runner profit-recalculation update to Sheets/Raiden_Profits.xlsx
# Updates currentAveragePrice / LevePrice / LeveProfit style columns (H:N) across several
# worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW) to match refreshed market-board figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 6
$ws.Range("H6").Value = 600
$ws.Range("I6").Value = 600
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1800
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1688

# Row 12
$ws.Range("H12").Value = 425.42856
$ws.Range("I12").Value = 395.6
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 395.6
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -225.6

# Row 15
$ws.Range("H15").Value = 855.65515
$ws.Range("I15").Value = 855.65515
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 2566.96545
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -2397.96545

# Row 38
$ws.Range("H38").Value = 396.16666
$ws.Range("I38").Value = 396.16666
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 1188.49998
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -816.4999800000001
$ws.Range("N38").ClearContents()

# Row 58
$ws.Range("H58").Value = 3198
$ws.Range("I58").Value = 3198
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 9594
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -9444
$ws.Range("N58").ClearContents()

# Row 87
$ws.Range("H87").Value = 27750
$ws.Range("I87").Value = 20000
$ws.Range("J87").Value = 30333.334
$ws.Range("K87").Value = 20000
$ws.Range("L87").Value = 30333.334
$ws.Range("M87").Value = -18752
$ws.Range("N87").Value = -32829.334

# Row 90
$ws.Range("H90").Value = 27750
$ws.Range("I90").Value = 20000
$ws.Range("J90").Value = 30333.334
$ws.Range("K90").Value = 60000
$ws.Range("L90").Value = 91000.00199999999
$ws.Range("M90").Value = -53760
$ws.Range("N90").Value = -103480.002

# Row 100
$ws.Range("H100").Value = 1236.75
$ws.Range("I100").Value = 1236.75
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1236.75
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -695.75

# Row 112
$ws.Range("H112").Value = 1435.2693
$ws.Range("I112").Value = 1213.3334
$ws.Range("J112").Value = 1737.909
$ws.Range("K112").Value = 3640.0002
$ws.Range("L112").Value = 5213.727000000001
$ws.Range("M112").Value = -2532.0002
$ws.Range("N112").Value = -7429.727000000001

# Row 125
$ws.Range("H125").Value = 2052.75
$ws.Range("I125").Value = 2052.75
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 18474.75
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -16014.75
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")

# Row 61
$ws.Range("H61").Value = 2621.0667
$ws.Range("I61").Value = 2545.64
$ws.Range("J61").Value = 2998.2
$ws.Range("K61").Value = 2545.64
$ws.Range("L61").Value = 2998.2
$ws.Range("M61").Value = -2333.64

# Row 102
$ws.Range("H102").Value = 2060.111
$ws.Range("I102").Value = 2017.625
$ws.Range("J102").Value = 2400
$ws.Range("K102").Value = 2017.625
$ws.Range("L102").Value = 2400
$ws.Range("M102").Value = -395.625

# Row 136
$ws.Range("H136").Value = 2621.0667
$ws.Range("I136").Value = 2545.64
$ws.Range("J136").Value = 2998.2
$ws.Range("K136").Value = 7636.92
$ws.Range("L136").Value = 8994.599999999999
$ws.Range("M136").Value = -5086.92

$ws = $wb.Worksheets.Item("BSM")

# Row 94
$ws.Range("H94").Value = 2276.0588
$ws.Range("I94").Value = 1076.9
$ws.Range("J94").Value = 3989.1428
$ws.Range("K94").Value = 1076.9
$ws.Range("L94").Value = 3989.1428
$ws.Range("M94").Value = -625.9000000000001

$ws = $wb.Worksheets.Item("CRP")

# Row 6
$ws.Range("H6").Value = 13368333
$ws.Range("I6").Value = 13368333
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 13368333
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -13368220

# Row 7
$ws.Range("H7").Value = 109.71429
$ws.Range("I7").Value = 53.6
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 53.6
$ws.Range("L7").Value = 250
$ws.Range("M7").Value = 59.4

# Row 17
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()

# Row 41
$ws.Range("H41").Value = 8571.286
$ws.Range("I41").Value = 5000
$ws.Range("J41").Value = 29999
$ws.Range("K41").Value = 5000
$ws.Range("L41").Value = 29999
$ws.Range("M41").Value = -4572
$ws.Range("N41").Value = -30855

# Row 51
$ws.Range("H51").Value = 40999.5
$ws.Range("I51").Value = 40999.5
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 40999.5
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -40263.5

# Row 61
$ws.Range("H61").Value = 40999.5
$ws.Range("I61").Value = 40999.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 40999.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -40651.5

$ws = $wb.Worksheets.Item("CUL")

# Row 9
$ws.Range("H9").Value = 17998.5
$ws.Range("I9").Value = 5998
$ws.Range("J9").Value = 29999
$ws.Range("K9").Value = 17994
$ws.Range("L9").Value = 89997
$ws.Range("M9").Value = -17770
$ws.Range("N9").Value = -90445

# Row 40
$ws.Range("H40").Value = 36
$ws.Range("I40").Value = 31.833334
$ws.Range("J40").Value = 52.666668
$ws.Range("K40").Value = 127.333336
$ws.Range("L40").Value = 210.666672
$ws.Range("M40").Value = -58.333336
$ws.Range("N40").Value = -348.666672

# Row 68
$ws.Range("H68").Value = 2157.5
$ws.Range("I68").Value = 1389.2
$ws.Range("J68").Value = 5999
$ws.Range("K68").Value = 4167.6
$ws.Range("L68").Value = 17997
$ws.Range("M68").Value = -3356.6
$ws.Range("N68").Value = -19619

# Row 69
$ws.Range("H69").Value = 15675.667
$ws.Range("I69").Value = 6999
$ws.Range("J69").Value = 20014
$ws.Range("K69").Value = 20997
$ws.Range("L69").Value = 60042
$ws.Range("M69").Value = -20186

# Row 71
$ws.Range("H71").Value = 2157.5
$ws.Range("I71").Value = 1389.2
$ws.Range("J71").Value = 5999
$ws.Range("K71").Value = 12502.8
$ws.Range("L71").Value = 53991
$ws.Range("M71").Value = -8446.800000000001
$ws.Range("N71").Value = -62103

# Row 72
$ws.Range("H72").Value = 15675.667
$ws.Range("I72").Value = 6999
$ws.Range("J72").Value = 20014
$ws.Range("K72").Value = 62991
$ws.Range("L72").Value = 180126
$ws.Range("M72").Value = -58935

$ws = $wb.Worksheets.Item("GSM")

# Row 39
$ws.Range("H39").Value = 53467.168
$ws.Range("I39").Value = 30237
$ws.Range("J39").Value = 58113.2
$ws.Range("K39").Value = 30237
$ws.Range("L39").Value = 58113.2
$ws.Range("M39").Value = -29705
$ws.Range("N39").Value = -59177.2

$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 1400
$ws.Range("I22").Value = 800
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 800
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -505
$ws.Range("N22").Value = -2590

# Row 27
$ws.Range("H27").Value = 1400
$ws.Range("I27").Value = 800
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 800
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -693
$ws.Range("N27").Value = -2214

# Row 60
$ws.Range("H60").Value = 155000
$ws.Range("I60").Value = 250000
$ws.Range("J60").Value = 60000
$ws.Range("K60").Value = 250000
$ws.Range("L60").Value = 60000
$ws.Range("M60").Value = -249491
$ws.Range("N60").Value = -61018

# Row 82
$ws.Range("H82").Value = 2845.75
$ws.Range("I82").Value = 2461
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 2461
$ws.Range("L82").Value = 4000
$ws.Range("M82").Value = -2100
$ws.Range("N82").Value = -4722

# Row 85
$ws.Range("H85").Value = 2845.75
$ws.Range("I85").Value = 2461
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2461
$ws.Range("L85").Value = 4000
$ws.Range("M85").Value = -1213
$ws.Range("N85").Value = -6496
